$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: remove bold/border/center style, clear A1 text
$ws.Range("A1:Y1").ClearFormats()
$ws.Range("A1").Value = ""

# Row 3
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 27
$ws.Range("P3").Value = 3
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 1
$ws.Range("W3").Value = 28

# Row 4
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 290
$ws.Range("P4").Value = 4
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 2
$ws.Range("W4").Value = 188

# Row 5
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 3295.06
$ws.Range("G5").Value = 233.62
$ws.Range("H5").Value = 867.5700000000001
$ws.Range("I5").Value = 233.62
$ws.Range("J5").Value = 533.89
$ws.Range("K5").Value = 533.89
$ws.Range("N5").Value = 333.68
$ws.Range("O5").Value = 74404.08
$ws.Range("P5").Value = 1334.75
$ws.Range("R5").Value = 233.62
$ws.Range("S5").Value = 233.62
$ws.Range("T5").Value = 867.5700000000001
$ws.Range("W5").Value = 60673.49

# Row 6
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 1.56
$ws.Range("G6").Value = 0.11
$ws.Range("H6").Value = 0.41
$ws.Range("I6").Value = 0.11
$ws.Range("J6").Value = 0.25
$ws.Range("K6").Value = 0.25
$ws.Range("N6").Value = 0.16
$ws.Range("O6").Value = 35.27
$ws.Range("P6").Value = 0.63
$ws.Range("R6").Value = 0.11
$ws.Range("S6").Value = 0.11
$ws.Range("T6").Value = 0.41
$ws.Range("W6").Value = 28.72

# Row 7
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 329.51
$ws.Range("G7").Value = 233.62
$ws.Range("H7").Value = 433.78
$ws.Range("I7").Value = 233.62
$ws.Range("J7").Value = 533.89
$ws.Range("K7").Value = 533.89
$ws.Range("N7").Value = 333.68
$ws.Range("O7").Value = 256.57
$ws.Range("P7").Value = 333.69
$ws.Range("R7").Value = 233.62
$ws.Range("S7").Value = 233.62
$ws.Range("T7").Value = 433.78
$ws.Range("W7").Value = 322.73

# Row 8
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 175.19
$ws.Range("G8").Value = 233.62
$ws.Range("H8").Value = 533.89
$ws.Range("I8").Value = 233.62
$ws.Range("J8").Value = 533.89
$ws.Range("K8").Value = 533.89
$ws.Range("N8").Value = 333.68
$ws.Range("O8").Value = 216.87
$ws.Range("P8").Value = 233.56
$ws.Range("R8").Value = 233.62
$ws.Range("S8").Value = 233.62
$ws.Range("T8").Value = 533.89
$ws.Range("W8").Value = 317
